$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.910.15"
$ws.Range("E2").Value = "  +5.47%  "
$ws.Range("D3").Value = "2.672.95"
$ws.Range("E3").Value = "  +6.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'610.51"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").Value = "'181.63"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("E9").Value = "  +17.43%  "
$ws.Range("D10").Value = "2.671.20"
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +5.41%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("E14").Value = "  +12.13%  "
$ws.Range("D15").Value = "3.121.77"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").Value = "'27.13"
$ws.Range("E16").Value = "  +5.32%  "
$ws.Range("D17").Value = "72.792.34"
$ws.Range("E17").Value = "  +5.61%  "
$ws.Range("D18").Value = "2.673.01"
$ws.Range("E18").Value = "  +6.76%  "
$ws.Range("D19").Value = "'385.92"
$ws.Range("E19").Value = "  +6.65%  "
$ws.Range("D20").Value = "'11.67"
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").Value = "'7.97"
$ws.Range("E21").Value = "  +5.89%  "
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("E23").Value = "  +24.37%  "
$ws.Range("D24").Value = "'73.55"
$ws.Range("E24").Value = "  +4.86%  "
$ws.Range("E25").Value = "  +7.68%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  +12.32%  "
$ws.Range("D28").Value = "2.810.93"
$ws.Range("E28").Value = "  +7.02%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "0.0₃0984"
$ws.Range("E30").Value = "  +11.86%  "
$ws.Range("D31").Value = "'548.91"
$ws.Range("E31").Value = "  +8.32%  "
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "  +5.98%  "
$ws.Range("E33").Value = "  +12.04%  "
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'164.85"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").Value = "'19.48"
$ws.Range("E37").Value = "  +4.26%  "
$ws.Range("E38").Value = "  +9.69%  "
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("D40").Value = "'19.16"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("E41").Value = "  +9.84%  "
$ws.Range("E42").Value = "  +8.60%  "
$ws.Range("E43").Value = "  +14.65%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +6.18%  "
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("D47").Value = "'153.87"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").Value = "'3.72"
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.549"
$ws.Range("E49").Value = "  +7.47%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0273"
$ws.Range("E50").Value = "  +12.30%  "
$ws.Range("D51").Value = "'1.73"
$ws.Range("E51").Value = "  +10.86%  "
